$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.765.97'
$ws.Range('E2').Value = '  -0.41%  '
$ws.Range('D3').Value = '1.633.61'
$ws.Range('E3').Value = '  -0.30%  '
$ws.Range('E4').Value = '  -0.20%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '215.04'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.29%  '
$ws.Range('E6').Value = '  -1.02%  '
$ws.Range('E7').Value = '  -0.31%  '
$ws.Range('E8').Value = '  -0.17%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.0637'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -1.10%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '19.69'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -3.69%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0787'
$ws.Range('D11').Style = "Normal"
$ws.Range('E12').Value = '  -0.35%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.640.94'
$ws.Range('E13').Value = '  +0.07%  '
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = '1.857.89'
$ws.Range('E14').Value = '  -0.33%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.559'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -0.44%  '
$ws.Range('D16').Value = '0.0₃0764'
$ws.Range('E16').Value = '  -0.63%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '62.77'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -0.97%  '
$ws.Range('D18').Value = '25.767.60'
$ws.Range('E18').Value = '  -0.40%  '
$ws.Range('E19').Value = '  -0.33%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '4.45'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +1.43%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '194.36'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +0.58%  '
$ws.Range('E22').Value = '  -0.58%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '6.27'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +1.74%  '
$ws.Range('E24').Value = '  -0.33%  '
$ws.Range('E25').Value = '  +4.11%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '142.66'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +2.37%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.123'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -0.35%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '6.88'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +0.42%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '15.55'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -0.31%  '
$ws.Range('E30').Value = '  -0.18%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.0494'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -0.62%  '
$ws.Range('E32').Value = '  +0.56%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.25'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -0.53%  '
$ws.Range('E34').Value = '  +0.26%  '
$ws.Range('E35').Value = '  -0.18%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.902'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -0.38%  '
$ws.Range('D37').Value = '1.129.51'
$ws.Range('E37').Value = '  -0.52%  '
$ws.Range('E38').Value = '  -1.61%  '
$ws.Range('E39').Value = '  -1.90%  '
$ws.Range('E40').Value = '  -1.05%  '
$ws.Range('E41').Value = '  +0.41%  '
$ws.Range('E42').Value = '  +1.91%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '99.99'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +0.74%  '
$ws.Range('E44').Value = '  +0.42%  '
$ws.Range('D45').Value = '1.767.63'
$ws.Range('E45').Value = '  -0.45%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '55.07'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -1.31%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.416'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -2.38%  '
$ws.Range('E49').Value = '  +0.01%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '7.58'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -3.22%  '
$ws.Range('B51').Value = 'SynthetixNetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.34'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +2.20%  '
